$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1 becomes "Grade", D1 is removed
$ws.Range("C1").Value = "Grade"
$ws.Range("D1").ClearContents()

# Populate row 2 with new data
$ws.Range("A2").Value = "Aarush"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "56"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "C1"

# Remove row 3 entirely (it previously held C3 = 370)
$ws.Range("A3:D3").ClearContents()
